$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Define the three new character styles (values taken from the target diff).
# Font.Color uses the standard OLE BGR packing (R + G*256 + B*65536), so
# 0x000080 (navy, R=00 G=00 B=80) is produced with decimal 8388608.
# ---------------------------------------------------------------------------
$ganStyle = $d.Styles.Add("GaNStyle", 2)
$ganStyle.Font.Name = "Calibri"
$ganStyle.Font.NameAscii = "Calibri"
$ganStyle.Font.Size = 14

$ganParagraph = $d.Styles.Add("GaNParagraph", 2)
$ganParagraph.Font.Name = "Calibri"
$ganParagraph.Font.NameAscii = "Calibri"
$ganParagraph.Font.Size = 10

$ganLinks = $d.Styles.Add("GaNLinks", 2)
$ganLinks.Font.Name = "Calibri"
$ganLinks.Font.NameAscii = "Calibri"
$ganLinks.Font.Bold = $true
$ganLinks.Font.Color = 8388608
$ganLinks.Font.Size = 9.5
$ganLinks.Font.Underline = 1

# ---------------------------------------------------------------------------
# Apply GaNStyle to every run containing the campaign-dates sentence
# (it appears four times in the document).
# ---------------------------------------------------------------------------
$datesText = "2022 Fechas de la campaña para constelación de orión: 16-25 de enero, 14-23 de febrero, 14-24 de marzo"
$rng = $d.Content
$found = $rng.Find.Execute($datesText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
while ($found) {
    $rng.Style = "GaNStyle"
    $rng.Collapse(0)
    $found = $rng.Find.Execute($datesText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
}

# ---------------------------------------------------------------------------
# Apply GaNParagraph to the introductory paragraph run.
# ---------------------------------------------------------------------------
$paragraphText = "Usted está participando en una campaña mundial para observar y registrar las estrellas visibles más débiles como un medio para medir la contaminación lumínica en un lugar determinado. Localizando y observando la  constelación de orión en el cielo nocturno y comparándolo con las cartas estelares, la gente de todo el mundo aprenderán cómo las luces de su comunidad contribuyen a la contaminación lumínica. Sus contribuciones a la base de datos en línea documentarán el cielo nocturno visible."
$rng2 = $d.Content
$found2 = $rng2.Find.Execute($paragraphText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $rng2.Style = "GaNParagraph"
}

# ---------------------------------------------------------------------------
# Apply GaNLinks to the map-link run.
# ---------------------------------------------------------------------------
$linkText = "(http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."
$rng3 = $d.Content
$found3 = $rng3.Find.Execute($linkText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $rng3.Style = "GaNLinks"
}
